$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210; this pushes the existing rows
# 210..231 down to 211..232 (shifting the whole data block, including the
# dimension, which grows from A1:R231 to A1:R232).
$ws.Rows.Item(210).Insert()

# Populate the newly-inserted row 210 with the new weekly price entry.
$ws.Range("A210").Value = 11
$ws.Range("B210").Value = "Vega Monumental Concepción"
$ws.Range("C210").Value = "Bíobío"
$ws.Range("D210").Value = 45142
$ws.Range("E210").Value = 8
$ws.Range("F210").Value = 100112043
$ws.Range("G210").Value = "Pepino ensalada"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 180
$ws.Range("K210").Value = 8000
$ws.Range("L210").Value = 9000
$ws.Range("M210").Value = 8333
$ws.Range("N210").Value = "$/caja 60 unidades"
$ws.Range("O210").Value = "Región de Arica y Parinacota"
$ws.Range("P210").Value = 139
$ws.Range("Q210").Value = 60
$ws.Range("R210").Value = "Hortaliza"
